$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-09-03 15:09:58"
$wsZhCn.Range("H2").Value = "2016-09-03 15:09:53"
$wsZhCn.Range("K2").Value = "2016-09-03 15:10:21"
$wsDeDe.Range("H2").Value = "2016-09-03 15:09:58"
$wsDeDe.Range("K2").Value = "2016-09-03 15:10:28"
